$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Drop the leading "S.No." column and the trailing "Status" column, pulling
# every other header one slot to the left (Admission No. now leads).
$ws.Range("A1").Value = "Admission No."
$ws.Range("B1").Value = "Student Name"
$ws.Range("C1").Value = "Parent / Guardian Name"
$ws.Range("D1").Value = "Class & Section"
$ws.Range("E1").Value = "Date Of Birth"
$ws.Range("F1").Value = "Blood Group"
$ws.Range("G1").Value = "Mobile No."
$ws.Range("H1").Value = "Alternative Mobile No."
$ws.Range("I1").Value = "RF ID Card No."
$ws.Range("J1").Value = "Address Line - 1"
$ws.Range("K1").Value = "Address Line - 2"
$ws.Range("L1").Value = "Address Line - 3"
$ws.Range("M1").Value = "City"
$ws.Range("N1").Value = "Pincode"

# The old O1 ("Pincode") / P1 ("Status") headers no longer exist.
$ws.Range("O1:P1").Clear()

# Widen column A (now "Admission No.") and split column C off from the old
# shared B:C width, matching what now lives there ("Parent / Guardian Name").
$ws.Columns("A").ColumnWidth = 17.7109375
$ws.Columns("C").ColumnWidth = 28.28515625

# The first data row now starts filling from A2 as well as B2; keep both
# formatted as text like the original B2 placeholder.
$ws.Range("A2:B2").NumberFormat = "@"

$ws.Range("A2").Select()
